$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.00033501180587336421
$ws.Range("A3").Value = 0.00018559227464720607
$ws.Range("H3").Value = 5.230769157409668
$ws.Range("A4").Value = 0.00014941954577807337
$ws.Range("H4").Value = 5.053255081176758
